# Reorders the worker/period detail rows (B16:G29) on "Hoja1" so that all
# rows for Periodo Mora "1707" come before the rows for Periodo Mora "1708",
# keeping the relative order of workers within each period group (the data
# for the 7 workers x 2 periods is unchanged; only the row order changes).
# This matches the "Actualiza base de datos EC y agrega parte 1 de nuevos
# estado de cuenta" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used in the detail table: B (Tipo Doc), C (N Doc), D (Nombre),
# E (Periodo Mora), F (Valor Mora), G (Salario Basico)
$cols = @("B", "C", "D", "E", "F", "G")

$firstRow = 16
$lastRow = 29

# Snapshot the current values of the 14 detail rows before writing anything
# back, using Value2 (reads/writes actual scalars reliably here).
$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    if ($rowVals["E"].ToString() -eq "1707") {
        $rowVals["SortKey"] = 0
    } else {
        $rowVals["SortKey"] = 1
    }
    $data += ,$rowVals
}

# Stable sort: "1707" rows first, then "1708" rows, preserving the original
# relative order of workers within each group.
$sorted = $data | Sort-Object -Property SortKey

# Write the reordered rows back.
$i = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = $sorted[$i]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $rowVals[$c]
    }
    $i = $i + 1
}
